$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Base")

# --- Update B3:B13 from numeric years to "Ano_YYYY" text labels ---
$ws.Range("B3").Value  = "Ano_2010"
$ws.Range("B4").Value  = "Ano_2011"
$ws.Range("B5").Value  = "Ano_2012"
$ws.Range("B6").Value  = "Ano_2013"
$ws.Range("B7").Value  = "Ano_2014"
$ws.Range("B8").Value  = "Ano_2015"
$ws.Range("B9").Value  = "Ano_2016"
$ws.Range("B10").Value = "Ano_2017"
$ws.Range("B11").Value = "Ano_2018"
$ws.Range("B12").Value = "Ano_2019"
$ws.Range("B13").Value = "Ano_2020"

# --- Define the named ranges (car-model columns + year rows) ---
$wb.Names.Add("Versa",    '=Base!$C$3:$C$13')
$wb.Names.Add("Kicks",    '=Base!$D$3:$D$13')
$wb.Names.Add("Corolla",  '=Base!$E$3:$E$13')
$wb.Names.Add("Etios",    '=Base!$F$3:$F$13')
$wb.Names.Add("Duster",   '=Base!$G$3:$G$13')
$wb.Names.Add("Logan",    '=Base!$H$3:$H$13')
$wb.Names.Add("Sandero",  '=Base!$I$3:$I$13')

$wb.Names.Add("Ano_2010", '=Base!$C$3:$I$3')
$wb.Names.Add("Ano_2011", '=Base!$C$4:$I$4')
$wb.Names.Add("Ano_2012", '=Base!$C$5:$I$5')
$wb.Names.Add("Ano_2013", '=Base!$C$6:$I$6')
$wb.Names.Add("Ano_2014", '=Base!$C$7:$I$7')
$wb.Names.Add("Ano_2015", '=Base!$C$8:$I$8')
$wb.Names.Add("Ano_2016", '=Base!$C$9:$I$9')
$wb.Names.Add("Ano_2017", '=Base!$C$10:$I$10')
$wb.Names.Add("Ano_2018", '=Base!$C$11:$I$11')
$wb.Names.Add("Ano_2019", '=Base!$C$12:$I$12')
$wb.Names.Add("Ano_2020", '=Base!$C$13:$I$13')

# --- K2 (brand selector) now points at Sandero instead of Etios ---
$ws.Range("K2").Value = "Sandero"

# --- K8 (year selector) now points at Ano_2020 instead of literal 2014 ---
$ws.Range("K8").Value = "Ano_2020"

# --- Brand block (rows 3-6) formulas driven off K2 ---
$ws.Range("L3").Formula = '=SUM(INDIRECT(K2))'
$ws.Range("L4").Formula = '=AVERAGE(INDIRECT(K2))'
$ws.Range("L5").Formula = '=MAX(INDIRECT(K2))'
$ws.Range("L6").Formula = '=MIN(INDIRECT(K2))'
$ws.Range("M5").Formula = '=INDEX(B3:B13, MATCH(L5, INDIRECT(K2), 0))'
$ws.Range("M6").Formula = '=INDEX(B3:B13, MATCH(L6, INDIRECT(K2),0))'

# M4 no longer shows the "-" placeholder now that it has a real sibling value
$ws.Range("M4").Value = ""

# --- Year block (rows 9-12) formulas driven off K8 ---
$ws.Range("L9").Formula  = '=SUM(INDIRECT(K8))'
$ws.Range("L10").Formula = '=AVERAGE(INDIRECT(K8))'
$ws.Range("L11").Formula = '=MAX(INDIRECT(K8))'
$ws.Range("L12").Formula = '=MIN(INDIRECT(K8))'
$ws.Range("M11").Formula = '=INDEX(C2:I2, MATCH(L11, INDIRECT(K8), 0))'
$ws.Range("M12").Formula = '=INDEX(C2:I2, MATCH(L12, INDIRECT(K8), 0))'

$wb.Application.Calculate()
